# Trade #75 closed at 2026-02-17 15:49:03 - unknown UNKNOWN +0.000%
#
# Updates:
#  - Summary sheet: roll the aggregate stats to reflect the newly closed trade.
#  - Strategy Status sheet: roll the MarketMaking strategy row stats.
#  - All Trades / MarketMaking sheets: append the new closed trade as row 76.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.11   # Current Capital
$summary.Range("B4").Value = 0.1       # Total P&L $
$summary.Range("B5").Value = 0.03      # Total P&L %
$summary.Range("B6").Value = 75        # Total Trades
$summary.Range("B8").Value = 40        # Losing Trades
$summary.Range("B9").Value = 32        # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet (row 4 = MarketMaking)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.11     # Capital
$status.Range("D4").Value = 75         # Trades
$status.Range("E4").Value = 0.1        # P&L $
$status.Range("F4").Value = 0.11       # P&L %
$status.Range("G4").Value = 32         # Win Rate %

# ---------------------------------------------------------------------------
# All Trades + MarketMaking sheets: append new closed trade as row 76
# ---------------------------------------------------------------------------
$tradeSheets = @($wb.Worksheets.Item("All Trades"), $wb.Worksheets.Item("MarketMaking"))

foreach ($ws in $tradeSheets) {
    $row = 76

    $ws.Cells.Item($row, 1).Value = 75                 # A: Trade #

    # Date / Time must stay plain text (not auto-parsed into date/time serials)
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"        # B: Date
    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = "15:48:57"          # C: Time

    $ws.Cells.Item($row, 4).Value = "MarketMaking"      # D: Strategy
    $ws.Cells.Item($row, 5).Value = "UP"                # E: Side
    $ws.Cells.Item($row, 6).Value = 0.71                # F: Entry Price
    $ws.Cells.Item($row, 7).Value = 0.67                # G: Exit Price
    $ws.Cells.Item($row, 8).Value = "CLOSED"            # H: Status
    $ws.Cells.Item($row, 9).Value = -5.6338             # I: P&L %
    $ws.Cells.Item($row, 10).Value = -0.04              # J: P&L $
    $ws.Cells.Item($row, 11).Value = 100.11             # K: Capital After
    $ws.Cells.Item($row, 12).Value = 0                  # L: Entry Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0                  # M: Exit Slippage (bps)
    $ws.Cells.Item($row, 14).Value = 0.6                # N: Confidence
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"  # O: Entry Reason
    $ws.Cells.Item($row, 16).Value = "early_exit"       # P: Exit Reason
    $ws.Cells.Item($row, 17).Value = 0.14                # Q: Duration (min)
}
